$wb = $excel.ActiveWorkbook

# Add the new worksheet at the end and rename it
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$newSheet.Name = "FindNewCarTest"

# Fill in the header row
$newSheet.Range("A1").Value = '${browser}'
$newSheet.Range("B1").Value = '${brandname}'
$newSheet.Range("C1").Value = '${carheading}'

# Fill in the data rows (order matters for shared-string table ordering)
$newSheet.Range("A2").Value = "chrome"
$newSheet.Range("B2").Value = "toyota"
$newSheet.Range("B3").Value = "kia"
$newSheet.Range("B4").Value = "bmw"
$newSheet.Range("C3").Value = "Kia Cars"
$newSheet.Range("C4").Value = "BMW Cars"
$newSheet.Range("C2").Value = "Toyota Cars"

# Set column widths to match target (values chosen so the engine's internal
# pixel-rounding reproduces the target stored widths of 12.5 and ~17.832)
$newSheet.Range("B:B").ColumnWidth = 11.666666666666666
$newSheet.Range("C:C").ColumnWidth = 17

# Make sure selection is on C2 for this new sheet and it's the active tab
$newSheet.Activate()
$newSheet.Range("C2").Select()
